$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.297455310821533
$ws.Range("B1").Value = 3.266421318054199
$ws.Range("C1").Value = 5.773363590240479
$ws.Range("D1").Value = 1.741952419281006
$ws.Range("E1").Value = 1.021064281463623
